# TC09_Bento_MultiFilter ... "startup" sheet fix:
#  - CasesTab / SamplesTab / FilesTab query cells (column B) each get an
#    "order By ... LIMIT 100" clause appended to their existing Cypher text.
#  - The sheet's active selection moves from B3 to C3 (and the prior
#    topLeftCell="A3" scroll-anchor goes away along with it).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- CasesTab row (row 2): append ORDER BY on study_subject_id ---
$ws.Range("B2").Value = $ws.Range("B2").Value2 + "`norder By ss.study_subject_id ASC LIMIT 100"

# --- SamplesTab row (row 3): append ORDER BY on sample_id (note the extra leading space,
#     exactly as authored in the source workbook) ---
$ws.Range("B3").Value = $ws.Range("B3").Value2 + "`n order By samp.sample_id ASC LIMIT 100"

# --- FilesTab row (row 4): append ORDER BY on file_name ---
$ws.Range("B4").Value = $ws.Range("B4").Value2 + "`norder By f.file_name ASC LIMIT 100"

# --- Update the active selection/view: was B3, now C3 ---
$ws.Range("C3").Select()
